$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the outlier-adjusted (MAD) imputation columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the existing header formatting (bold, centered, bordered) from E1
# onto the new header cells so they reuse the same cell style.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill the new data columns with FALSE boolean values for rows 2-4
$ws.Range("F2:H4").Value = $false
